$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Step5" spreadsheet step row (whole row), shifting everything below up by one.
$ws.Rows(24).Delete()

# After the shift, the test table (rows 27-29) only keeps a single "Result" column (G),
# referencing Step6 instead of Step5; clear the now-redundant second result column (H).
$ws.Range("G27").Value = "_res_.`$Step6"
$ws.Range("H27").ClearContents()

$ws.Range("H28").ClearContents()

$ws.Range("G29").Value = "java.lang.Object"
$ws.Range("H29").ClearContents()

# Add explanatory comment describing why the expected result is Object (not Object[]).
$ws.Range("J29").Value = "// true ? Policy[] : Double -> should be Object (not Object[]). Please, don't do this"
